$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 2 and 3: "Profitable..." now comes before "Total Net Worth..."
$ws.Range("A2").Value = "Profitable in latest Financial Statement (Latest Full Year)"
$ws.Range("B2").Value = "Yes"

$ws.Range("A3").Value = "Total Net Worth (Total Equity)"
$ws.Range("B3").Value = "Negative"

# Update ratio values
$ws.Range("B4").Value = "0.34"
$ws.Range("B5").Value = "-1.53"

# Widen the default column width for the sheet
$ws.Cells.EntireColumn.AutoFit() | Out-Null
$ws.Application.ActiveWindow | Out-Null
$ws.DefaultColumnWidth = 50.21625
